$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 130
$ws.Range("J5").Value = 200
$ws.Range("L5").Value = 200
$ws.Range("N5").Value = -430
$ws.Range("H8").Value = 2135.125
$ws.Range("I8").Value = 2791
$ws.Range("J8").Value = 167.5
$ws.Range("K8").Value = 8373
$ws.Range("L8").Value = 502.5
$ws.Range("M8").Value = -8234
$ws.Range("N8").Value = -780.5
$ws.Range("H9").Value = 247.5
$ws.Range("I9").Value = 221.75
$ws.Range("K9").Value = 221.75
$ws.Range("M9").Value = -52.75
$ws.Range("H11").Value = 1185.9375
$ws.Range("I11").Value = 1185.9375
$ws.Range("K11").Value = 1185.9375
$ws.Range("M11").Value = -1045.9375
$ws.Range("H64").Value = 5079.6
$ws.Range("I64").Value = 5348.5
$ws.Range("K64").Value = 5348.5
$ws.Range("M64").Value = -5100.5
$ws.Range("H67").Value = 5079.6
$ws.Range("I67").Value = 5348.5
$ws.Range("K67").Value = 5348.5
$ws.Range("M67").Value = -4490.5
$ws.Range("H70").Value = 5690.636
$ws.Range("I70").Value = 6999.25
$ws.Range("K70").Value = 20997.75
$ws.Range("M70").Value = -20727.75
$ws.Range("H73").Value = 5690.636
$ws.Range("I73").Value = 6999.25
$ws.Range("K73").Value = 20997.75
$ws.Range("M73").Value = -20061.75
$ws.Range("H101").Value = 2244
$ws.Range("I101").Value = 2492.5
$ws.Range("K101").Value = 7477.5
$ws.Range("M101").Value = -5855.5
$ws.Range("H105").Value = 100997.25
$ws.Range("J105").Value = 100997.25
$ws.Range("L105").Value = 100997.25
$ws.Range("N105").Value = -107985.25
$ws.Range("H107").Value = 1762.6111
$ws.Range("J107").Value = 2413.1667
$ws.Range("L107").Value = 2413.1667
$ws.Range("N107").Value = -6253.1667
$ws.Range("H137").Value = 4374.048
$ws.Range("I137").Value = 1710.7
$ws.Range("J137").Value = 6795.273
$ws.Range("K137").Value = 5132.1
$ws.Range("L137").Value = 20385.819
$ws.Range("M137").Value = -2582.1
$ws.Range("N137").Value = -25485.819
$ws.Range("H138").Value = 2939.92
$ws.Range("I138").Value = 1047.25
$ws.Range("J138").Value = 3104.5
$ws.Range("K138").Value = 3141.75
$ws.Range("L138").Value = 9313.5
$ws.Range("M138").Value = 1998.25
$ws.Range("N138").Value = -19593.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H105").Value = 50000
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H125").Value = 82940.8
$ws.Range("J125").Value = 82940.8
$ws.Range("L125").Value = 82940.8
$ws.Range("N125").Value = -92780.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 110000
$ws.Range("J108").Value = 110000
$ws.Range("L108").Value = 110000
$ws.Range("N108").Value = -117680
$ws.Range("H110").Value = 129000
$ws.Range("J110").Value = 129000
$ws.Range("L110").Value = 129000
$ws.Range("N110").Value = -137180

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 5000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 25000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -31240
$ws.Range("H97").Value = 103000
$ws.Range("J97").Value = 103000
$ws.Range("L97").Value = 103000
$ws.Range("N97").Value = -104982
$ws.Range("H134").Value = 459027.47
$ws.Range("I134").Value = 626862.2
$ws.Range("J134").Value = 11468.167
$ws.Range("K134").Value = 1880586.6
$ws.Range("L134").Value = 34404.501
$ws.Range("M134").Value = -1878051.6
$ws.Range("N134").Value = -39474.501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 84.75
$ws.Range("I7").Value = 84.75
$ws.Range("K7").Value = 254.25
$ws.Range("M7").Value = -142.25
$ws.Range("H34").Value = 474.66666
$ws.Range("I34").Value = 474.66666
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1423.99998
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1339.99998
$ws.Range("N34").ClearContents()
$ws.Range("H37").Value = 80998.664
$ws.Range("J37").Value = 80998.664
$ws.Range("L37").Value = 242995.992
$ws.Range("N37").Value = -243219.992
$ws.Range("H55").Value = 500
$ws.Range("I55").Value = 500
$ws.Range("K55").Value = 1500
$ws.Range("M55").Value = -1323
$ws.Range("H80").Value = 4049.0908
$ws.Range("I80").Value = 2996.5
$ws.Range("K80").Value = 8989.5
$ws.Range("M80").Value = -8053.5
$ws.Range("H83").Value = 4049.0908
$ws.Range("I83").Value = 2996.5
$ws.Range("K83").Value = 26968.5
$ws.Range("M83").Value = -22288.5
$ws.Range("H104").Value = 4622.75
$ws.Range("I104").Value = 4497
$ws.Range("K104").Value = 13491
$ws.Range("M104").Value = -10870
$ws.Range("H107").Value = 540.75
$ws.Range("I107").Value = 285
$ws.Range("K107").Value = 855
$ws.Range("M107").Value = 1065
$ws.Range("H108").Value = 1442.3334
$ws.Range("I108").Value = 1413.5
$ws.Range("J108").Value = 1500
$ws.Range("K108").Value = 4240.5
$ws.Range("L108").Value = 4500
$ws.Range("M108").Value = -1360.5
$ws.Range("N108").Value = -10260
$ws.Range("H131").Value = 8507.333000000001
$ws.Range("I131").Value = 9114.666999999999
$ws.Range("J131").Value = 7900
$ws.Range("K131").Value = 27344.001
$ws.Range("L131").Value = 23700
$ws.Range("M131").Value = -22304.001
$ws.Range("N131").Value = -33780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 248.13333
$ws.Range("I2").Value = 40.125
$ws.Range("J2").Value = 485.85715
$ws.Range("K2").Value = 40.125
$ws.Range("L2").Value = 485.85715
$ws.Range("M2").Value = 72.875
$ws.Range("N2").Value = -711.85715
$ws.Range("H139").Value = 75000
$ws.Range("I139").Value = 97000
$ws.Range("J139").Value = 64000
$ws.Range("K139").Value = 97000
$ws.Range("L139").Value = 64000
$ws.Range("M139").Value = -91860
$ws.Range("N139").Value = -74280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 98333.336
$ws.Range("J87").Value = 95000
$ws.Range("L87").Value = 95000
$ws.Range("N87").Value = -97246
$ws.Range("H90").Value = 98333.336
$ws.Range("J90").Value = 95000
$ws.Range("L90").Value = 285000
$ws.Range("N90").Value = -296232
$ws.Range("H104").Value = 43145.168
$ws.Range("J104").Value = 43145.168
$ws.Range("L104").Value = 43145.168
$ws.Range("N104").Value = -50133.168
$ws.Range("H122").Value = 4748.1562
$ws.Range("I122").Value = 3979.2964
$ws.Range("K122").Value = 11937.8892
$ws.Range("M122").Value = -9487.889200000001
$ws.Range("H132").Value = 525474.25
$ws.Range("I132").Value = 557169.7
$ws.Range("J132").Value = 335301.66
$ws.Range("K132").Value = 1671509.1
$ws.Range("L132").Value = 1005904.98
$ws.Range("M132").Value = -1668979.1
$ws.Range("N132").Value = -1010964.98
$ws.Range("H136").Value = 71846.06
$ws.Range("J136").Value = 167287
$ws.Range("L136").Value = 501861
$ws.Range("N136").Value = -506961

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H62").Value = 25006848
$ws.Range("I62").Value = 7830.5
$ws.Range("J62").Value = 50005864
$ws.Range("K62").Value = 7830.5
$ws.Range("L62").Value = 50005864
$ws.Range("M62").Value = -7206.5
$ws.Range("N62").Value = -50007112
$ws.Range("H65").Value = 25006848
$ws.Range("I65").Value = 7830.5
$ws.Range("J65").Value = 50005864
$ws.Range("K65").Value = 39152.5
$ws.Range("L65").Value = 250029320
$ws.Range("M65").Value = -36032.5
$ws.Range("N65").Value = -250035560
$ws.Range("H92").Value = 67030
$ws.Range("J92").Value = 67030
$ws.Range("L92").Value = 67030
$ws.Range("N92").Value = -72022
$ws.Range("H100").Value = 1349.6
$ws.Range("I100").Value = 1349.6
$ws.Range("K100").Value = 2699.2
$ws.Range("M100").Value = -2158.2
$ws.Range("H113").Value = 898
$ws.Range("I113").Value = 926.9091
$ws.Range("K113").Value = 2780.7273
$ws.Range("M113").Value = -610.7273
$ws.Range("H122").Value = 4603.077
$ws.Range("I122").Value = 2566
$ws.Range("K122").Value = 7698
$ws.Range("M122").Value = -5248
$ws.Range("H136").Value = 1569.5
$ws.Range("I136").Value = 1650.1428
$ws.Range("J136").Value = 1005
$ws.Range("K136").Value = 4950.428400000001
$ws.Range("L136").Value = 3015
$ws.Range("M136").Value = -2400.428400000001
$ws.Range("N136").Value = -8115
